# Logic tree input file updated.
#
# Two new "Possible_Problem" header rows are inserted into the logic tree,
# mirroring the existing row 4 pattern (Node1 = the relevant Problem node,
# Relationship = "Possible_Problem", Node2 = the generic/default
# Possible_Problem answer), for each of the two remaining Problem nodes that
# did not yet have one:
#   - "Does the engine run too hot? ..." (new row 7, pushes old rows 7-11 down to 8-12)
#   - "Are you seeing the coolant on the ground every time you park? ..." (new row 12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new row, above the existing "Does the engine run too
# hot?" block (old row 7) ---
$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = "Problem:Does the engine run too hot? (check temperature gauge on the dashboard) (Please answer as: Yes, No, Not Sure, Fluctuates)"
$ws.Range("B7").Value = "Possible_Problem"
$ws.Range("C7").Value = "Possible_Problem:35% Radiator (Hoses)`n30% Water Pump`n20% Thermostat Housing`n15% Cylinder Head Gasket"
$ws.Range("C7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 230.4

# --- Insert the second new row, above the existing "Are you seeing the
# coolant..." block (old row 11, now row 12 after the first insert) ---
$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = "Problem:Are you seeing the coolant on the ground every time you park? (Please answer as: Yes, No)"
$ws.Range("B12").Value = "Possible_Problem"
$ws.Range("C12").Value = "Possible_Problem:35% Radiator (Hoses)`n30% Water Pump`n20% Thermostat Housing`n15% Cylinder Head Gasket"
$ws.Range("C12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 230.4

# --- Update the view: selection / top-left cell moved to reflect the new
# bottom of the sheet ---
$ws.Range("C12").Select()
$excel.ActiveWindow.ScrollRow = 14
